$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: A4 792 -> 805, B4 "No results" -> "NO Tname and no Bname"
$ws.Range("A4").Value = 805
$ws.Range("B4").Value = "NO Tname and no Bname"

# Row 6: A6 1258 -> 965, B6 "Genetically Modified Declaration Code" -> "NO Tname and no Bname"
$ws.Range("A6").Value = 965
$ws.Range("B6").Value = "NO Tname and no Bname"
